# Update the dSF (column F) values per the repulled/pushed data and mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = -6
$ws.Range("F6").Value = 0
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = 9
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = 0
